$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 22: B22 was stored as text "19"; convert it to a true numeric value 19
$ws.Cells.Item(22, 2).Value = 19

# New row 23 with student data
$ws.Cells.Item(23, 1).Value = "Vanda"
$ws.Cells.Item(23, 2).NumberFormat = "@"
$ws.Cells.Item(23, 2).Value = "19"
$ws.Cells.Item(23, 3).Value = "Male"
$ws.Cells.Item(23, 4).Value = "Phnom Penh"
$ws.Cells.Item(23, 5).Value = "Singer"
$ws.Cells.Item(23, 6).Value = "image\bdef9cf1bf584c58ac5ec5e2ff915481.png"
